$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix window type length test data bugs
# Row 3: Window # (column B) was 2, should be 1
$ws.Range("B3").Value = 1

# Row 5: Start Time (column D) corrected
$ws.Range("D5").Value = 55450 / 86400

# Update the active selection on the sheet to D6
$ws.Range("D6").Select()
